$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the reporting period / timestamp headers ---
$ws.Range("J7").Value = "2024-11-01 - 2024-11-08"
$ws.Range("E8").Value = "2024-11-08T14:58:"

# --- Row 16 (first ledger line) ---
$ws.Range("C16").Value = "USD"
$ws.Range("G16").Value = "4800000.0"
$ws.Range("N16").Value = "2"

# --- Row 18 (second ledger line) ---
$ws.Range("C18").Value = "SP"
$ws.Range("G18").Value = "4000000.0"
$ws.Range("P18").Value = "عبد الجواد اللايح"
$ws.Range("U18").Value = "251716"
$ws.Range("X18").Value = "339"

# --- Row 20 (third ledger line) ---
$ws.Range("C20").Value = "USD"
$ws.Range("G20").Value = "880000.0"
$ws.Range("P20").Value = "عاصم شيخو"
$ws.Range("U20").Value = "829983"
$ws.Range("X20").Value = "51"

# --- Drop the old fourth ledger line (row 22) and the trailing spacer row
#     (row 23); row 21 (the blank spacer row between line 3 and the removed
#     line 4) absorbs the freed-up height. ---
$ws.Rows.Item(22).Delete()
$ws.Rows.Item(22).Delete()
$ws.Rows.Item(21).RowHeight = 802
